$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.581.90'
$ws.Range("E2").Value = '  -0.73%  '

$ws.Range("D3").Value = '2.291.72'
$ws.Range("E3").Value = '  -1.74%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '495.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.28%  '

$ws.Range("D9").Value = '2.292.52'
$ws.Range("E9").Value = '  -1.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0953'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.77%  '

$ws.Range("E11").Value = '  +2.04%  '

$ws.Range("E12").Value = '  +2.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.48%  '

$ws.Range("D14").Value = '2.695.09'
$ws.Range("E14").Value = '  -1.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.80%  '

$ws.Range("D16").Value = '54.501.80'
$ws.Range("E16").Value = '  -0.81%  '

$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").Value = '2.290.24'
$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("E19").Value = '  +3.03%  '

$ws.Range("E20").Value = '  +2.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '305.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("E24").Value = '  -2.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.57%  '

$ws.Range("E26").Value = '  +1.11%  '

$ws.Range("E27").Value = '  +5.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.373'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '

$ws.Range("D29").Value = '2.371.13'
$ws.Range("E29").Value = '  -2.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("E33").Value = '  -2.27%  '

$ws.Range("E34").Value = '  +2.73%  '

$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.19%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.60'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("E39").Value = '  +2.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.860'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.11%  '

$ws.Range("E41").Value = '  +0.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.39%  '

$ws.Range("E43").Value = '  +1.74%  '

$ws.Range("E44").Value = '  +1.52%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '128.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.77%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.85%  '

$ws.Range("E48").Value = '  +0.72%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.549'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '242.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '

$ws.Range("E51").Value = '  +1.51%  '
